# Update the representative-coordinates label and the Site 2 lat/lon values
# on Sheet1 (wind and solar capacity factor corrections for Site 2 / TX).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = "34.223, -102.245"
$ws.Range("C4").Value = 34.334000000000003
$ws.Range("C5").Value = -102.245

[void]$ws.Range("C6").Select()
